$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "26.060.78"
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "  +0.70%  "
$r.Style = "Normal"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.644.81"
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "  +0.96%  "
$r.Style = "Normal"
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "  +1.20%  "
$r.Style = "Normal"
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "216.75"
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "  +1.14%  "
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "  +1.03%  "
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "  +0.99%  "
$r.Style = "Normal"
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "  +0.55%  "
$r.Style = "Normal"
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "  +1.37%  "
$r.Style = "Normal"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "  +0.27%  "
$r.Style = "Normal"
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0797"
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "  +0.95%  "
$r.Style = "Normal"
$r = $ws.Range("B12")
$r.NumberFormat = "@"
$r.Value = "WrappedliquidstakedEther2.0"
$r.Style = "Normal"
$r = $ws.Range("C12")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$r.Style = "Normal"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "1.873.29"
$r.Style = "Normal"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "  +0.99%  "
$r.Style = "Normal"
$r = $ws.Range("B13")
$r.NumberFormat = "@"
$r.Value = "Polkadot"
$r.Style = "Normal"
$r = $ws.Range("C13")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$r.Style = "Normal"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "4.29"
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "  +1.48%  "
$r.Style = "Normal"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "1.666.48"
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "  +1.38%  "
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "  +0.32%  "
$r.Style = "Normal"
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.0₃0765"
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "  +1.23%  "
$r.Style = "Normal"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "63.33"
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "  +0.96%  "
$r.Style = "Normal"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "26.138.50"
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "  +1.03%  "
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "  +1.06%  "
$r.Style = "Normal"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "193.17"
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "  +0.03%  "
$r.Style = "Normal"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "4.34"
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "  -0.86%  "
$r.Style = "Normal"
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "9.94"
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "  +0.00%  "
$r.Style = "Normal"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "  -0.28%  "
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "  +1.30%  "
$r.Style = "Normal"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "144.45"
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "  +1.62%  "
$r.Style = "Normal"
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "  +1.19%  "
$r.Style = "Normal"
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "  +3.95%  "
$r.Style = "Normal"
$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "  +0.85%  "
$r.Style = "Normal"
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "15.53"
$r.Style = "Normal"
$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "  +0.48%  "
$r.Style = "Normal"
$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "  +1.24%  "
$r.Style = "Normal"
$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = "  -0.44%  "
$r.Style = "Normal"
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.30"
$r.Style = "Normal"
$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "  -0.44%  "
$r.Style = "Normal"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "3.27"
$r.Style = "Normal"
$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "  +1.51%  "
$r.Style = "Normal"
$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "  +2.20%  "
$r.Style = "Normal"
$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "  -3.22%  "
$r.Style = "Normal"
$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "  +0.60%  "
$r.Style = "Normal"
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "1.133.75"
$r.Style = "Normal"
$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = "  +0.22%  "
$r.Style = "Normal"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.540"
$r.Style = "Normal"
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "  -1.72%  "
$r.Style = "Normal"
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "  +0.13%  "
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "  +0.49%  "
$r.Style = "Normal"
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "  +1.01%  "
$r.Style = "Normal"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "99.52"
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "  +0.43%  "
$r.Style = "Normal"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.798"
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "  -0.77%  "
$r.Style = "Normal"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "1.781.54"
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "  +0.97%  "
$r.Style = "Normal"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.0₆0116"
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "  +4.34%  "
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "56.56"
$r.Style = "Normal"
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "  +0.87%  "
$r.Style = "Normal"
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "  +1.09%  "
$r.Style = "Normal"
$r = $ws.Range("B48")
$r.NumberFormat = "@"
$r.Value = "RenderToken"
$r.Style = "Normal"
$r = $ws.Range("C48")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$r.Style = "Normal"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "1.45"
$r.Style = "Normal"
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "  -0.05%  "
$r.Style = "Normal"
$r = $ws.Range("B49")
$r.NumberFormat = "@"
$r.Value = "EnergySwap"
$r.Style = "Normal"
$r = $ws.Range("C49")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$r.Style = "Normal"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "7.73"
$r.Style = "Normal"
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "  +1.26%  "
$r.Style = "Normal"
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "  +0.79%  "
$r.Style = "Normal"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0957"
$r.Style = "Normal"
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "  -0.39%  "
$r.Style = "Normal"
